$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 481; this shifts rows 481:521 down to 482:522
# and preserves formatting already present (e.g. the D-column date style).
$ws.Rows(481).Insert()

# Populate the newly inserted row 481 with the new weekly data point.
$ws.Cells.Item(481, 1).Value = 5
$ws.Cells.Item(481, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(481, 3).Value = "Maule"
$ws.Cells.Item(481, 4).Value = 45223
$ws.Cells.Item(481, 5).Value = 7
$ws.Cells.Item(481, 6).Value = 100112008
$ws.Cells.Item(481, 7).Value = "Coliflor"
$ws.Cells.Item(481, 8).Value = "Sin especificar"
$ws.Cells.Item(481, 9).Value = "Primera"
$ws.Cells.Item(481, 10).Value = 5000
$ws.Cells.Item(481, 11).Value = 600
$ws.Cells.Item(481, 12).Value = 600
$ws.Cells.Item(481, 13).Value = 600
$ws.Cells.Item(481, 14).Value = "$/unidad"
$ws.Cells.Item(481, 15).Value = "Región del Maule"
$ws.Cells.Item(481, 16).Value = 600
$ws.Cells.Item(481, 17).Value = 1
$ws.Cells.Item(481, 18).Value = "Hortaliza"
